# Apply the workbook edits described by the diff:
#  - Selenium sheet (data table + becomes the active/selected tab, selection -> B7)
#  - Main sheet loses the "active tab" status, selection -> B19
#  - Temperature sheet selection -> H20
#  - Workbook first visible scrolled tab -> Sediment (index 9, best-effort)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the Selenium sheet's data table (rows 2-7, columns C/D/E, plus the
#    newly populated B7:E7 row).
# ---------------------------------------------------------------------------
$wsSel = $wb.Worksheets.Item("Selenium")

$wsSel.Range("C2").Value = 0
$wsSel.Range("D2").Value = 100
$wsSel.Range("E2").Value = 100

$wsSel.Range("C3").Value = 5
$wsSel.Range("D3").Value = 50
$wsSel.Range("E3").Value = 100

$wsSel.Range("C4").Value = 10
$wsSel.Range("D4").Value = 25
$wsSel.Range("E4").Value = 75

$wsSel.Range("C5").Value = 5
$wsSel.Range("D5").Value = 0
$wsSel.Range("E5").Value = 50

$wsSel.Range("C6").Value = 0
$wsSel.Range("D6").Value = 0
$wsSel.Range("E6").Value = 0

$wsSel.Range("B7").Value = 0
$wsSel.Range("C7").Value = 0
$wsSel.Range("D7").Value = 0
$wsSel.Range("E7").Value = 0

# ---------------------------------------------------------------------------
# 2. Update the view/selection state of the various sheets. Excel only keeps
#    one sheet "tabSelected" at a time, and a sheet's Range.Select() only
#    "sticks" for the sheet that is active when it's called - so visit each
#    sheet, make it active, and set its selection; the LAST sheet activated
#    ends up as the workbook's active tab (matches activeTab="13" / Selenium
#    tabSelected="1" in the target).
# ---------------------------------------------------------------------------

# Main: keep its existing data, just move the selection and make sure it is
# no longer the active tab (handled automatically once another sheet is
# activated below).
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Activate()
$wsMain.Range("B19").Select()

# Temperature: move the selection to H20.
$wsTemp = $wb.Worksheets.Item("Temperature")
$wsTemp.Activate()
$wsTemp.Range("H20").Select()

# Scroll the workbook tab strip so Sediment is the first displayed tab
# (firstSheet="9" in the target workbookView).
$wsSediment = $wb.Worksheets.Item("Sediment")
$win = $excel.ActiveWindow
$win.ScrollWorkbookTabs(9, 0)

# Selenium becomes the active/selected sheet, with the selection on B7.
$wsSel.Activate()
$wsSel.Range("B7").Select()

Write-Output "Edit applied"
